# "Tambah Peserta" — add a running-number ("No") column to the MENDAFTAR sheet.
#
# Before: columns A:E = NIM, Nama, E-mail, CP, Tanggal Seleksi I (table "Table2").
# After:  a new column A = "No" (1..25) is inserted in front, pushing the
#         existing table to B:F, and the new "No" column becomes its own
#         one-column table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MENDAFTAR")

# Insert a new blank column before column A; this shifts the existing
# NIM/Nama/E-mail/CP/Tanggal data (and the Table2 list object) from A:E to B:F.
$ws.Columns("A:A").Insert() | Out-Null

# The engine doesn't auto-resize the list object's range on column insert,
# so move Table2 onto its new B4:F29 footprint explicitly.
$tbl = $ws.ListObjects.Item("Table2")
$tbl.Resize($ws.Range("B4:F29")) | Out-Null

# Header + sequential numbering for the new "No" column (rows 5-29 match the
# 25 data rows of Table2).
$ws.Range("A4").Value = "No"
for ($i = 1; $i -le 25; $i++) {
    $ws.Cells.Item(4 + $i, 1).Value = $i
}

# Turn the new column into its own table, matching the light style used by
# the workbook's other tables.
$noTable = $ws.ListObjects.Add(1, $ws.Range("A4:A29"), 0, 1)
$noTable.TableStyle = "TableStyleLight1"

# The two e-mail hyperlinks lived in column C (Tanggal... no, E-mail) and are
# now one column over, in D; re-anchor them there.
$ws.Hyperlinks.Delete() | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:obbie.christian@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D9"), "mailto:joeladlino26@yahoo.com") | Out-Null
# Re-adding a hyperlink resets the cell's font size to the workbook default;
# put it back to match the sheet's normal 10pt text.
$ws.Range("D5").Font.Size = 10
$ws.Range("D9").Font.Size = 10

# Move the active selection the way the author left it.
$ws.Range("B13").Select() | Out-Null
